$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("AB1").Value = "13-10-2020"
$ws.Range("AB1").Font.Bold = $true
$ws.Range("AB1").HorizontalAlignment = -4108
$ws.Range("AB1").VerticalAlignment = -4160
$ws.Range("AB1").Borders.LineStyle = 1
$ws.Range("AB2").Value = 55
$ws.Range("AB3").Value = 6256
$ws.Range("AB4").Value = 24
$ws.Range("AB5").Value = 826
$ws.Range("AB6").Value = 955
$ws.Range("AB7").Value = 192
$ws.Range("AB8").Value = 1286
$ws.Range("AB9").Value = 2
$ws.Range("AB10").Value = 5809
$ws.Range("AB11").Value = 511
$ws.Range("AB12").Value = 3574
$ws.Range("AB13").Value = 1592
$ws.Range("AB14").Value = 251
$ws.Range("AB15").Value = 1333
$ws.Range("AB16").Value = 798
$ws.Range("AB17").Value = 10036
$ws.Range("AB18").Value = 1025
$ws.Range("AB19").Value = 64
$ws.Range("AB20").Value = 2645
$ws.Range("AB21").Value = 40514
$ws.Range("AB22").Value = 93
$ws.Range("AB23").Value = 64
$ws.Range("AB24").Value = 0
$ws.Range("AB25").Value = 18
$ws.Range("AB26").Value = 1040
$ws.Range("AB27").Value = 565
$ws.Range("AB28").Value = 3860
$ws.Range("AB29").Value = 1665
$ws.Range("AB30").Value = 57
$ws.Range("AB31").Value = 10314
$ws.Range("AB32").Value = 1233
$ws.Range("AB33").Value = 317
$ws.Range("AB34").Value = 762
$ws.Range("AB35").Value = 6438
$ws.Range("AB36").Value = 5682
